$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("alpha2F")

# Append a new row (row 16), copying the formatting of the row above (row 15)
# so the new row matches the existing style (bold/bordered/centered "A" column style).
$ws.Range("A15:M15").Copy()
$ws.Range("A16:M16").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16").Value = 1.13768898836368
$ws.Range("D16").Value = 1.344620042544543
$ws.Range("E16").Value = 0.8372164374120621
$ws.Range("F16").Value = 1.13768898836368
$ws.Range("G16").Value = 1.08840550983146
$ws.Range("H16").Value = 0.976955512593477
$ws.Range("I16").Value = 0.9024898308117507
$ws.Range("J16").Value = 1.344620042544543
$ws.Range("K16").Value = 1.090918239978302
$ws.Range("L16").Value = 1.114303614170991
$ws.Range("M16").Value = 1.047896053592829
